$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells with new computed values
$ws.Range("D2").Value = 0.0009573850256213977
$ws.Range("E2").Value = 0.04776818181139628

$ws.Range("D4").Value = 0.01062241905879754

$ws.Range("D5").Value = 0.009789598803911358

# Add new row 6 for DWA - copy formatting from row 2's A cell (style index 1)
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A6").Value = "DWA"
$ws.Range("B6").Value = 75
$ws.Range("C6").Value = 0.7488421052631579
$ws.Range("D6").Value = 0.02732875324688148
$ws.Range("E6").Value = 0.7187556818398444
